# The workbook originally has two sheets:
#   sheet1.xml (rId1) name="hotel_info"  (sheetId=1) - 9 cols x 2 rows (incl header)
#   sheet2.xml (rId2) name="review_info" (sheetId=2) - 25 cols x 1 row (header only)
#
# Target state:
#   sheet1.xml (rId1) name="review_info" (sheetId=1) - 25 cols x 1 row (header only)
#   sheet2.xml (rId2) name="hotel_info"  (sheetId=2) - 10 cols x 2 rows (incl header),
#       with a new "State" column inserted between Hotel_Name and City, value "Louisiana"
#
# We keep the physical sheet1.xml/sheet2.xml (and thus rId/sheetId) slots fixed and
# instead swap the cell data + names between them, so sheetId stays aligned with r:id
# exactly like the target diff shows.

$wb = $excel.ActiveWorkbook
$sheetA = $wb.Worksheets.Item("hotel_info")   # physically sheet1.xml / rId1
$sheetB = $wb.Worksheets.Item("review_info")  # physically sheet2.xml / rId2

# --- capture current hotel_info data from sheetA (9 cols x 2 rows) ---
$hotelHeaders = @()
for ($c = 1; $c -le 9; $c++) {
    $hotelHeaders += ,$sheetA.Cells.Item(1, $c).Value()
}
$hotelRow2 = @()
for ($c = 1; $c -le 9; $c++) {
    $hotelRow2 += ,$sheetA.Cells.Item(2, $c).Value()
}

# --- capture current review_info headers from sheetB (25 cols x 1 row) ---
$reviewHeaders = @()
for ($c = 1; $c -le 25; $c++) {
    $reviewHeaders += ,$sheetB.Cells.Item(1, $c).Value()
}

# --- clear both sheets entirely ---
$sheetA.Cells.Clear() | Out-Null
$sheetB.Cells.Clear() | Out-Null

# rename out of the way first to avoid name collisions while swapping
$sheetA.Name = "__tmp_sheetA__"
$sheetB.Name = "__tmp_sheetB__"

# --- write review_info headers into sheetA, rename to review_info ---
for ($c = 1; $c -le 25; $c++) {
    $sheetA.Cells.Item(1, $c).Value = $reviewHeaders[$c - 1]
}
$sheetA.Name = "review_info"

# --- write hotel_info data (with new State column) into sheetB, rename to hotel_info ---
# New column order: STR, Hotel_Name, State, City, Zip, TA_ReviewURL,
#                    Tripadvisor_Hotel_Name, English_Reviews_num, Local_Rank, Total_Reviews_num
$newHotelHeaders = @($hotelHeaders[0], $hotelHeaders[1], "State", $hotelHeaders[2], $hotelHeaders[3], $hotelHeaders[4], $hotelHeaders[5], $hotelHeaders[6], $hotelHeaders[7], $hotelHeaders[8])
$newHotelRow2    = @($hotelRow2[0],    $hotelRow2[1],    "Louisiana", $hotelRow2[2], $hotelRow2[3], $hotelRow2[4], $hotelRow2[5], $hotelRow2[6], $hotelRow2[7], $hotelRow2[8])

for ($c = 1; $c -le 10; $c++) {
    $sheetB.Cells.Item(1, $c).Value = $newHotelHeaders[$c - 1]
}

# Columns H, I, J (English_Reviews_num, Local_Rank, Total_Reviews_num) hold
# digit-only strings in the source data ("1155", "59", "1185") - force Text
# number format first so COM doesn't silently coerce them back to numbers.
$sheetB.Cells.Item(2, 8).NumberFormat = "@"
$sheetB.Cells.Item(2, 9).NumberFormat = "@"
$sheetB.Cells.Item(2, 10).NumberFormat = "@"
for ($c = 1; $c -le 10; $c++) {
    $sheetB.Cells.Item(2, $c).Value = $newHotelRow2[$c - 1]
}
$sheetB.Name = "hotel_info"

# Column A (STR / date-of-scraping-like serial) and column E (Zip) are numeric
# in the source data - re-assert that now that the whole row has been written,
# since COM may have stringified them when they passed through the generic
# object array above.
$sheetB.Cells.Item(2, 1).NumberFormat = "General"
$sheetB.Cells.Item(2, 5).NumberFormat = "General"
$sheetB.Cells.Item(2, 1).Value = [double]$hotelRow2[0]
$sheetB.Cells.Item(2, 5).Value = [double]$hotelRow2[3]
